$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.737.87'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '2.714.28'
$ws.Range('E3').Value = '  +2.01%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '114.69'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '329.85'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.529'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.559'
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.50'
$ws.Range('E10').Value = '  -2.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.22'
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0824'
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.129'
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.64'
$ws.Range('E14').Value = '  +3.42%  '
$ws.Range('D15').Value = '3.129.26'
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').Value = '2.672.94'
$ws.Range('E16').Value = '  +1.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.877'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').Value = '50.527.69'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.81'
$ws.Range('E19').Value = '  +3.75%  '
$ws.Range('B20').Value = 'ImmutableX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.96'
$ws.Range('E20').Value = '  +0.62%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.81'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').Value = '0.0₃0959'
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '276.16'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.04'
$ws.Range('E24').Value = '  -3.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.59'
$ws.Range('E25').Value = '  -0.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.64'
$ws.Range('E26').Value = '  -1.29%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.32'
$ws.Range('E28').Value = '  +2.78%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.141'
$ws.Range('E30').Value = '  -1.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.55'
$ws.Range('E31').Value = '  -3.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.86'
$ws.Range('E32').Value = '  -0.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.58'
$ws.Range('E33').Value = '  +1.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0825'
$ws.Range('E34').Value = '  +1.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '19.40'
$ws.Range('E35').Value = '  -1.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.997'
$ws.Range('E36').Value = '  -0.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.01'
$ws.Range('E37').Value = '  -0.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.08'
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.20'
$ws.Range('E39').Value = '  +1.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.75'
$ws.Range('E40').Value = '  +5.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '128.70'
$ws.Range('E41').Value = '  +3.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0351'
$ws.Range('E42').Value = '  +9.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.32'
$ws.Range('E43').Value = '  +4.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.113'
$ws.Range('E44').Value = '  -0.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.41'
$ws.Range('E45').Value = '  +2.05%  '
$ws.Range('D46').Value = '2.094.53'
$ws.Range('E46').Value = '  -0.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.16'
$ws.Range('E47').Value = '  +8.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.05'
$ws.Range('E49').Value = '  -0.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.40'
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '59.83'
$ws.Range('E51').Value = '  -0.50%  '
